$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.805.03'
$ws.Range("E2").Value = '  +1.96%  '
$ws.Range("D3").Value = '3.024.30'
$ws.Range("E3").Value = '  +0.70%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '510.91'
$ws.Range("E5").Value = '  +0.12%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.36'
$ws.Range("E6").Value = '  +0.52%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  +0.67%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.60'
$ws.Range("E9").Value = '  +0.35%  '
$ws.Range("E10").Value = '  +1.51%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.367'
$ws.Range("E11").Value = '  +3.39%  '
$ws.Range("D12").Value = '3.536.39'
$ws.Range("E12").Value = '  +0.62%  '
$ws.Range("E13").Value = '  +0.71%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.67'
$ws.Range("E14").Value = '  +3.81%  '
$ws.Range("E15").Value = '  +5.58%  '
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '57.769.43'
$ws.Range("E16").Value = '  +1.82%  '
$ws.Range("B17").Value = 'Polkadot'
$ws.Range("C17").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.29'
$ws.Range("E17").Value = '  +5.90%  '
$ws.Range("D18").Value = '3.024.86'
$ws.Range("E18").Value = '  +0.84%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.93'
$ws.Range("E19").Value = '  +3.41%  '
$ws.Range("E20").Value = '  +1.83%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '330.64'
$ws.Range("E21").Value = '  +0.23%  '
$ws.Range("E22").Value = '  +0.17%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.72'
$ws.Range("E23").Value = '  -0.78%  '
$ws.Range("E24").Value = '  +3.51%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '64.87'
$ws.Range("E25").Value = '  +3.18%  '
$ws.Range("E26").Value = '  -2.93%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  -0.51%  '
$ws.Range("D28").Value = '0.0₃0925'
$ws.Range("E28").Value = '  +1.28%  '
$ws.Range("E29").Value = '  +2.16%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.37'
$ws.Range("E30").Value = '  +3.97%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.82'
$ws.Range("E31").Value = '  +1.48%  '
$ws.Range("E32").Value = '  -5.07%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.71'
$ws.Range("E33").Value = '  +0.20%  '
$ws.Range("E34").Value = '  +3.94%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '154.78'
$ws.Range("E35").Value = '  -0.30%  '
$ws.Range("E36").Value = '  +4.09%  '
$ws.Range("E37").Value = '  +0.58%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '24.63'
$ws.Range("E38").Value = '  +1.73%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0680'
$ws.Range("E39").Value = '  -0.76%  '
$ws.Range("D40").Value = '3.056.16'
$ws.Range("E40").Value = '  +0.75%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '37.71'
$ws.Range("E41").Value = '  +2.03%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.87'
$ws.Range("E42").Value = '  +5.69%  '
$ws.Range("E43").Value = '  -0.08%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.653'
$ws.Range("E44").Value = '  +0.68%  '
$ws.Range("E45").Value = '  +0.09%  '
$ws.Range("D46").Value = '2.233.35'
$ws.Range("E46").Value = '  -1.48%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.989'
$ws.Range("E47").Value = '  -1.37%  '
$ws.Range("E48").Value = '  +4.16%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0240'
$ws.Range("E49").Value = '  +0.16%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '19.62'
$ws.Range("E50").Value = '  +1.27%  '
$ws.Range("E51").Value = '  -5.64%  '
